$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.966.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "'1.634.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'211.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'23.49"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").Value = "'0.0615"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "'1.864.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "'1.633.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "'65.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'27.962.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'232.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "'7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  -6.17%  "
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'155.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "'6.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").Value = "'0.0482"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "'3.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").Value = "'3.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'1.410.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "'1.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.05%  "
$ws.Range("D37").Value = "'2.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("D39").Value = "'0.557"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'66.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").Value = "'5.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "'1.775.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").Value = "'88.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").Value = "'0.0₆0104"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").Value = "'0.0998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("E51").Value = "  -0.38%  "
